$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-07 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-08 Friday", 2) | Out-Null
$d.Content.Find.Execute("91-19=72", $true, $false, $false, $false, $false, $true, 1, $false, "22+7=29", 2) | Out-Null
$d.Content.Find.Execute("38-1=37", $true, $false, $false, $false, $false, $true, 1, $false, "96-57=39", 2) | Out-Null
$d.Content.Find.Execute("86-82=4", $true, $false, $false, $false, $false, $true, 1, $false, "7+15=22", 2) | Out-Null
$d.Content.Find.Execute("45-25=20", $true, $false, $false, $false, $false, $true, 1, $false, "65+8=73", 2) | Out-Null
$d.Content.Find.Execute("24+32=56", $true, $false, $false, $false, $false, $true, 1, $false, "85-59=26", 2) | Out-Null
$d.Content.Find.Execute("54+3=57", $true, $false, $false, $false, $false, $true, 1, $false, "76-27=49", 2) | Out-Null
$d.Content.Find.Execute("28+7=35", $true, $false, $false, $false, $false, $true, 1, $false, "90-50=40", 2) | Out-Null
$d.Content.Find.Execute("25+70=95", $true, $false, $false, $false, $false, $true, 1, $false, "85-9=76", 2) | Out-Null
$d.Content.Find.Execute("23+24=47", $true, $false, $false, $false, $false, $true, 1, $false, "81-3=78", 2) | Out-Null
$d.Content.Find.Execute("40-16=24", $true, $false, $false, $false, $false, $true, 1, $false, "39+59=98", 2) | Out-Null
$d.Content.Find.Execute("16+2=18", $true, $false, $false, $false, $false, $true, 1, $false, "65-38=27", 2) | Out-Null
$d.Content.Find.Execute("83-15=68", $true, $false, $false, $false, $false, $true, 1, $false, "87+3=90", 2) | Out-Null
$d.Content.Find.Execute("32+62=94", $true, $false, $false, $false, $false, $true, 1, $false, "10+50=60", 2) | Out-Null
$d.Content.Find.Execute("26+21=47", $true, $false, $false, $false, $false, $true, 1, $false, "99-63=36", 2) | Out-Null
$d.Content.Find.Execute("78-5=73", $true, $false, $false, $false, $false, $true, 1, $false, "67-35=32", 2) | Out-Null
$d.Content.Find.Execute("29-27=2", $true, $false, $false, $false, $false, $true, 1, $false, "22+28=50", 2) | Out-Null
$d.Content.Find.Execute("97-8=89", $true, $false, $false, $false, $false, $true, 1, $false, "97-49=48", 2) | Out-Null
$d.Content.Find.Execute("48+22=70", $true, $false, $false, $false, $false, $true, 1, $false, "19+58=77", 2) | Out-Null
$d.Content.Find.Execute("96-2=94", $true, $false, $false, $false, $false, $true, 1, $false, "76-47=29", 2) | Out-Null
$d.Content.Find.Execute("46-8=38", $true, $false, $false, $false, $false, $true, 1, $false, "50-5=45", 2) | Out-Null
$d.Content.Find.Execute("34-29=5", $true, $false, $false, $false, $false, $true, 1, $false, "31-1=30", 2) | Out-Null
$d.Content.Find.Execute("59-28=31", $true, $false, $false, $false, $false, $true, 1, $false, "29+35=64", 2) | Out-Null
$d.Content.Find.Execute("84-27=57", $true, $false, $false, $false, $false, $true, 1, $false, "38-36=2", 2) | Out-Null
$d.Content.Find.Execute("8+52=60", $true, $false, $false, $false, $false, $true, 1, $false, "64-17=47", 2) | Out-Null
$d.Content.Find.Execute("6+88=94", $true, $false, $false, $false, $false, $true, 1, $false, "78+21=99", 2) | Out-Null
$d.Content.Find.Execute("22-3=19", $true, $false, $false, $false, $false, $true, 1, $false, "12+31=43", 2) | Out-Null
$d.Content.Find.Execute("2+57=59", $true, $false, $false, $false, $false, $true, 1, $false, "26+32=58", 2) | Out-Null
$d.Content.Find.Execute("8+30=38", $true, $false, $false, $false, $false, $true, 1, $false, "69-14=55", 2) | Out-Null
$d.Content.Find.Execute("87-13=74", $true, $false, $false, $false, $false, $true, 1, $false, "61-53=8", 2) | Out-Null
$d.Content.Find.Execute("59+20=79", $true, $false, $false, $false, $false, $true, 1, $false, "90+0=90", 2) | Out-Null
$d.Content.Find.Execute("62-33=29", $true, $false, $false, $false, $false, $true, 1, $false, "97-35=62", 2) | Out-Null
$d.Content.Find.Execute("28+27=55", $true, $false, $false, $false, $false, $true, 1, $false, "5+42=47", 2) | Out-Null
$d.Content.Find.Execute("1+34=35", $true, $false, $false, $false, $false, $true, 1, $false, "69-40=29", 2) | Out-Null
$d.Content.Find.Execute("80-78=2", $true, $false, $false, $false, $false, $true, 1, $false, "13+46=59", 2) | Out-Null
$d.Content.Find.Execute("5+79=84", $true, $false, $false, $false, $false, $true, 1, $false, "92-13=79", 2) | Out-Null
$d.Content.Find.Execute("29+30=59", $true, $false, $false, $false, $false, $true, 1, $false, "40-9=31", 2) | Out-Null
$d.Content.Find.Execute("45+11=56", $true, $false, $false, $false, $false, $true, 1, $false, "95-75=20", 2) | Out-Null
$d.Content.Find.Execute("50-39=11", $true, $false, $false, $false, $false, $true, 1, $false, "93-26=67", 2) | Out-Null
$d.Content.Find.Execute("62+23=85", $true, $false, $false, $false, $false, $true, 1, $false, "47+31=78", 2) | Out-Null
$d.Content.Find.Execute("66-33=33", $true, $false, $false, $false, $false, $true, 1, $false, "36-5=31", 2) | Out-Null
$d.Content.Find.Execute("86-7=79", $true, $false, $false, $false, $false, $true, 1, $false, "51+1=52", 2) | Out-Null
$d.Content.Find.Execute("74-1=73", $true, $false, $false, $false, $false, $true, 1, $false, "75-20=55", 2) | Out-Null
$d.Content.Find.Execute("40-31=9", $true, $false, $false, $false, $false, $true, 1, $false, "97-50=47", 2) | Out-Null
$d.Content.Find.Execute("91-47=44", $true, $false, $false, $false, $false, $true, 1, $false, "85-57=28", 2) | Out-Null
$d.Content.Find.Execute("79-70=9", $true, $false, $false, $false, $false, $true, 1, $false, "77+9=86", 2) | Out-Null
$d.Content.Find.Execute("27+22=49", $true, $false, $false, $false, $false, $true, 1, $false, "42-28=14", 2) | Out-Null
$d.Content.Find.Execute("58+32=90", $true, $false, $false, $false, $false, $true, 1, $false, "26-13=13", 2) | Out-Null
$d.Content.Find.Execute("28+46=74", $true, $false, $false, $false, $false, $true, 1, $false, "49+32=81", 2) | Out-Null
$d.Content.Find.Execute("46+44=90", $true, $false, $false, $false, $false, $true, 1, $false, "2+76=78", 2) | Out-Null
$d.Content.Find.Execute("73-26=47", $true, $false, $false, $false, $false, $true, 1, $false, "15+81=96", 2) | Out-Null
$d.Content.Find.Execute("49+16=65", $true, $false, $false, $false, $false, $true, 1, $false, "62-57=5", 2) | Out-Null
$d.Content.Find.Execute("27+58=85", $true, $false, $false, $false, $false, $true, 1, $false, "87-87=0", 2) | Out-Null
$d.Content.Find.Execute("32+28=60", $true, $false, $false, $false, $false, $true, 1, $false, "23-18=5", 2) | Out-Null
$d.Content.Find.Execute("36-32=4", $true, $false, $false, $false, $false, $true, 1, $false, "82+4=86", 2) | Out-Null
$d.Content.Find.Execute("51-6=45", $true, $false, $false, $false, $false, $true, 1, $false, "95-52=43", 2) | Out-Null
$d.Content.Find.Execute("63-3=60", $true, $false, $false, $false, $false, $true, 1, $false, "52-23=29", 2) | Out-Null
$d.Content.Find.Execute("57+23=80", $true, $false, $false, $false, $false, $true, 1, $false, "36+33=69", 2) | Out-Null
$d.Content.Find.Execute("92-79=13", $true, $false, $false, $false, $false, $true, 1, $false, "97-58=39", 2) | Out-Null
$d.Content.Find.Execute("87-30=57", $true, $false, $false, $false, $false, $true, 1, $false, "82-62=20", 2) | Out-Null
$d.Content.Find.Execute("56+5=61", $true, $false, $false, $false, $false, $true, 1, $false, "44-15=29", 2) | Out-Null
$d.Content.Find.Execute("72-27=45", $true, $false, $false, $false, $false, $true, 1, $false, "41-8=33", 2) | Out-Null
$d.Content.Find.Execute("68-62=6", $true, $false, $false, $false, $false, $true, 1, $false, "43+32=75", 2) | Out-Null
$d.Content.Find.Execute("49+26=75", $true, $false, $false, $false, $false, $true, 1, $false, "32+0=32", 2) | Out-Null
$d.Content.Find.Execute("34+59=93", $true, $false, $false, $false, $false, $true, 1, $false, "90-88=2", 2) | Out-Null
$d.Content.Find.Execute("74-8=66", $true, $false, $false, $false, $false, $true, 1, $false, "1+41=42", 2) | Out-Null
$d.Content.Find.Execute("85+13=98", $true, $false, $false, $false, $false, $true, 1, $false, "27+34=61", 2) | Out-Null
$d.Content.Find.Execute("49-40=9", $true, $false, $false, $false, $false, $true, 1, $false, "58-0=58", 2) | Out-Null
$d.Content.Find.Execute("7+4=11", $true, $false, $false, $false, $false, $true, 1, $false, "90-14=76", 2) | Out-Null
$d.Content.Find.Execute("16+12=28", $true, $false, $false, $false, $false, $true, 1, $false, "30-22=8", 2) | Out-Null
$d.Content.Find.Execute("23+61=84", $true, $false, $false, $false, $false, $true, 1, $false, "91-58=33", 2) | Out-Null
$d.Content.Find.Execute("44-19=25", $true, $false, $false, $false, $false, $true, 1, $false, "16+16=32", 2) | Out-Null
$d.Content.Find.Execute("28+2=30", $true, $false, $false, $false, $false, $true, 1, $false, "24-16=8", 2) | Out-Null
$d.Content.Find.Execute("74-19=55", $true, $false, $false, $false, $false, $true, 1, $false, "55-24=31", 2) | Out-Null
$d.Content.Find.Execute("95-50=45", $true, $false, $false, $false, $false, $true, 1, $false, "31-17=14", 2) | Out-Null
$d.Content.Find.Execute("70-23=47", $true, $false, $false, $false, $false, $true, 1, $false, "89-32=57", 2) | Out-Null
$d.Content.Find.Execute("7-4=3", $true, $false, $false, $false, $false, $true, 1, $false, "45+44=89", 2) | Out-Null
$d.Content.Find.Execute("88-79=9", $true, $false, $false, $false, $false, $true, 1, $false, "18+33=51", 2) | Out-Null
$d.Content.Find.Execute("10+33=43", $true, $false, $false, $false, $false, $true, 1, $false, "65+14=79", 2) | Out-Null
$d.Content.Find.Execute("53-4=49", $true, $false, $false, $false, $false, $true, 1, $false, "22+73=95", 2) | Out-Null
$d.Content.Find.Execute("34+60=94", $true, $false, $false, $false, $false, $true, 1, $false, "90-64=26", 2) | Out-Null
$d.Content.Find.Execute("38+24=62", $true, $false, $false, $false, $false, $true, 1, $false, "26+12=38", 2) | Out-Null
$d.Content.Find.Execute("85-78=7", $true, $false, $false, $false, $false, $true, 1, $false, "45+33=78", 2) | Out-Null
$d.Content.Find.Execute("4+48=52", $true, $false, $false, $false, $false, $true, 1, $false, "77-65=12", 2) | Out-Null
$d.Content.Find.Execute("56-10=46", $true, $false, $false, $false, $false, $true, 1, $false, "8+87=95", 2) | Out-Null
$d.Content.Find.Execute("10+27=37", $true, $false, $false, $false, $false, $true, 1, $false, "36-30=6", 2) | Out-Null
$d.Content.Find.Execute("72-48=24", $true, $false, $false, $false, $false, $true, 1, $false, "57-49=8", 2) | Out-Null
$d.Content.Find.Execute("24+66=90", $true, $false, $false, $false, $false, $true, 1, $false, "2+28=30", 2) | Out-Null
$d.Content.Find.Execute("45+15=60", $true, $false, $false, $false, $false, $true, 1, $false, "75-8=67", 2) | Out-Null
$d.Content.Find.Execute("65-24=41", $true, $false, $false, $false, $false, $true, 1, $false, "90-81=9", 2) | Out-Null
$d.Content.Find.Execute("78-33=45", $true, $false, $false, $false, $false, $true, 1, $false, "49+3=52", 2) | Out-Null
$d.Content.Find.Execute("33+47=80", $true, $false, $false, $false, $false, $true, 1, $false, "22-6=16", 2) | Out-Null
$d.Content.Find.Execute("67-11=56", $true, $false, $false, $false, $false, $true, 1, $false, "35+14=49", 2) | Out-Null
$d.Content.Find.Execute("4+76=80", $true, $false, $false, $false, $false, $true, 1, $false, "34-13=21", 2) | Out-Null
$d.Content.Find.Execute("40+47=87", $true, $false, $false, $false, $false, $true, 1, $false, "0+38=38", 2) | Out-Null
$d.Content.Find.Execute("2+55=57", $true, $false, $false, $false, $false, $true, 1, $false, "19-1=18", 2) | Out-Null
$d.Content.Find.Execute("59-42=17", $true, $false, $false, $false, $false, $true, 1, $false, "46-26=20", 2) | Out-Null
$d.Content.Find.Execute("65-58=7", $true, $false, $false, $false, $false, $true, 1, $false, "2+44=46", 2) | Out-Null
$d.Content.Find.Execute("65-46=19", $true, $false, $false, $false, $false, $true, 1, $false, "99-83=16", 2) | Out-Null
$d.Content.Find.Execute("16+80=96", $true, $false, $false, $false, $false, $true, 1, $false, "66+13=79", 2) | Out-Null
$d.Content.Find.Execute("66-55=11", $true, $false, $false, $false, $false, $true, 1, $false, "10+46=56", 2) | Out-Null
